$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") rows 2-127 currently hold 43; update to 243
$ws.Range("B2:B127").Value = 243
